# Generate Report for Handoff
# Refresh the "Overview" / "zh-cn" / "de-de" localization-status report with
# two additional source files (a markdown file and its two PNG dependents),
# updating the existing handoff row and appending rows 3-4 on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Shared data for the three "source files" this handoff run covers.
# ---------------------------------------------------------------------------
$pngA      = "04869b9e-4d17-4fb4-9288-6c4fb763e8cf.png"
$mdFile    = "06878a13-d130-4c4f-964c-168968a99dc8.md"
$pngB      = "37c74c71-09b8-442f-a1f1-75a2bc8fc461.png"

$status       = "Ready for handoff"
$overviewDate = "2016-44-12 10:44:08"

$zhHandoffFile  = "e2a21ff923e4e631faf124413bb63357ebc0330c.png"
$mdZhHandoffFile = "06878a13-d130-4c4f-964c-168968a99dc8.4bcc1834b3db8645c6f846e3a8733d046296a8c3.zh-cn.xlf"
$pngBZhHandoffFile = "c802d86fa7b58755480afab6e96d1cd8bd0b5bcb.png"
$zhDate = "2016-03-12 10:44:05"

$deHandoffFile  = "e2a21ff923e4e631faf124413bb63357ebc0330c.png"
$mdDeHandoffFile = "06878a13-d130-4c4f-964c-168968a99dc8.4bcc1834b3db8645c6f846e3a8733d046296a8c3.de-de.xlf"
$pngBDeHandoffFile = "c802d86fa7b58755480afab6e96d1cd8bd0b5bcb.png"
$deDate = "2016-03-12 10:44:08"

$epoch = "0001-01-01 00:00:00"
$dependencyFrom = "e2e\06878a13-d130-4c4f-964c-168968a99dc8.md"

$repoBlob   = "https://github.com/OpenLocalizationTest/oltest/blob/a7a4f8c2989548e208d46047de7b4589a407c89c/e2e"
$zhHandoffBlob = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c41e93776d241924123b29314623ae2de2eabc47/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deHandoffBlob = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1cd76ff7a85cd11c7994459cefe75487484e9be0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

function Set-CellHyperlink {
    param($ws, $cell, [string]$target, [string]$display)

    $range = $ws.Range($cell)
    if ($range.Hyperlinks.Count -gt 0) {
        $range.Hyperlinks.Delete()
    }
    $ws.Hyperlinks.Add($range, $target, "", "", $display) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = $pngA
$ov.Range("B2").Value = $status
$ov.Range("C2").Value = $status
$ov.Range("D2").Value = $overviewDate

$ov.Range("A3").Value = $mdFile
$ov.Range("B3").Value = $status
$ov.Range("C3").Value = $status
$ov.Range("D3").Value = $overviewDate

$ov.Range("A4").Value = $pngB
$ov.Range("B4").Value = $status
$ov.Range("C4").Value = $status
$ov.Range("D4").Value = $overviewDate

Set-CellHyperlink $ov "A2" "$repoBlob/$pngA" $pngA
Set-CellHyperlink $ov "A3" "$repoBlob/$mdFile" $mdFile
Set-CellHyperlink $ov "A4" "$repoBlob/$pngB" $pngB

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = $pngA
$zh.Range("B2").Value = ".png"
$zh.Range("C2").Value = $status
$zh.Range("D2").Value = $zhHandoffFile
$zh.Range("E2").Value = $zhDate
$zh.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H2").Value = $epoch
$zh.Range("I2").Value = "IsDependency"
$zh.Range("J2").Value = $dependencyFrom

$zh.Range("A3").Value = $mdFile
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = $status
$zh.Range("D3").Value = $mdZhHandoffFile
$zh.Range("E3").Value = $zhDate
$zh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H3").Value = $epoch
$zh.Range("I3").Value = "Include"

$zh.Range("A4").Value = $pngB
$zh.Range("B4").Value = ".png"
$zh.Range("C4").Value = $status
$zh.Range("D4").Value = $pngBZhHandoffFile
$zh.Range("E4").Value = $zhDate
$zh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H4").Value = $epoch
$zh.Range("I4").Value = "IsDependency"
$zh.Range("J4").Value = $dependencyFrom

Set-CellHyperlink $zh "A2" "$repoBlob/$pngA" $pngA
Set-CellHyperlink $zh "B2" "$repoBlob/$pngA" ".png"
Set-CellHyperlink $zh "D2" "$zhHandoffBlob/$zhHandoffFile" $zhHandoffFile

Set-CellHyperlink $zh "A3" "$repoBlob/$mdFile" $mdFile
Set-CellHyperlink $zh "B3" "$repoBlob/$mdFile" ".md"
Set-CellHyperlink $zh "D3" "$zhHandoffBlob/$mdZhHandoffFile" $mdZhHandoffFile

Set-CellHyperlink $zh "A4" "$repoBlob/$pngB" $pngB
Set-CellHyperlink $zh "B4" "$repoBlob/$pngB" ".png"
Set-CellHyperlink $zh "D4" "$zhHandoffBlob/$pngBZhHandoffFile" $pngBZhHandoffFile

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = $pngA
$de.Range("B2").Value = ".png"
$de.Range("C2").Value = $status
$de.Range("D2").Value = $deHandoffFile
$de.Range("E2").Value = $deDate
$de.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H2").Value = $epoch
$de.Range("I2").Value = "IsDependency"
$de.Range("J2").Value = $dependencyFrom

$de.Range("A3").Value = $mdFile
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = $status
$de.Range("D3").Value = $mdDeHandoffFile
$de.Range("E3").Value = $deDate
$de.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H3").Value = $epoch
$de.Range("I3").Value = "Include"

$de.Range("A4").Value = $pngB
$de.Range("B4").Value = ".png"
$de.Range("C4").Value = $status
$de.Range("D4").Value = $pngBDeHandoffFile
$de.Range("E4").Value = $deDate
$de.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H4").Value = $epoch
$de.Range("I4").Value = "IsDependency"
$de.Range("J4").Value = $dependencyFrom

Set-CellHyperlink $de "A2" "$repoBlob/$pngA" $pngA
Set-CellHyperlink $de "B2" "$repoBlob/$pngA" ".png"
Set-CellHyperlink $de "D2" "$deHandoffBlob/$deHandoffFile" $deHandoffFile

Set-CellHyperlink $de "A3" "$repoBlob/$mdFile" $mdFile
Set-CellHyperlink $de "B3" "$repoBlob/$mdFile" ".md"
Set-CellHyperlink $de "D3" "$deHandoffBlob/$mdDeHandoffFile" $mdDeHandoffFile

Set-CellHyperlink $de "A4" "$repoBlob/$pngB" $pngB
Set-CellHyperlink $de "B4" "$repoBlob/$pngB" ".png"
Set-CellHyperlink $de "D4" "$deHandoffBlob/$pngBDeHandoffFile" $pngBDeHandoffFile
